$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All B:E data cells in this sheet are plain text (the source xml uses
# t="inlineStr" throughout), even cells whose text happens to look like a
# plain number (e.g. "218.82"). Excel's COM .Value setter auto-converts
# numeric-looking strings into real numbers, which would change the cell's
# stored type/representation. Prefix those with an apostrophe (the classic
# force-text input trick) so the value is kept as text, then reset the style
# back to Normal so the cell doesn't pick up a stray quote-prefix format --
# none of these data cells carried an explicit style originally.

$ws.Cells.Item(2, 4).Value = '26.863.32'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +0.51%  '
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(3, 4).Value = '1.642.99'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +0.10%  '
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(5, 4).Value = '''218.82'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.09%  '
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(6, 4).Value = '''0.497'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.48%  '
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.16%  '
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.11%  '
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -1.02%  '
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(10, 4).Value = '''19.26'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +0.71%  '
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(11, 4).Value = '''0.0845'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.50%  '
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(12, 4).Value = '1.871.18'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +0.01%  '
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(13, 4).Value = '1.643.98'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.04%  '
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.20%  '
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +0.23%  '
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(16, 4).Value = '''65.45'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +1.85%  '
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(17, 4).Value = '26.868.45'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.49%  '
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(18, 4).Value = '0.0₃0733'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(19, 4).Value = '''215.25'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.63%  '
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.03%  '
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(21, 4).Value = '''6.67'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +6.59%  '
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(22, 4).Value = '''4.36'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.12%  '
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -1.58%  '
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(24, 4).Value = '''9.21'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -1.22%  '
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(25, 4).Value = '''147.79'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +2.05%  '
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.12%  '
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.18%  '
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(28, 4).Value = '''7.24'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +2.11%  '
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(29, 4).Value = '''15.72'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +0.74%  '
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -0.04%  '
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(31, 4).Value = '''1.20'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +1.71%  '
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +1.81%  '
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +0.10%  '
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(34, 4).Value = '1.282.68'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -0.73%  '
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +0.88%  '
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +0.17%  '
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(37, 4).Value = '''0.0173'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +0.07%  '
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -0.17%  '
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.53%  '
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.06%  '
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +0.08%  '
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +0.20%  '
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(43, 4).Value = '1.782.79'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.57%  '
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -5.98%  '
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(45, 4).Value = '''92.74'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +1.61%  '
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -0.22%  '
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -0.07%  '
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.28%  '
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 2).Style = "Normal"
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49, 3).Style = "Normal"
$ws.Cells.Item(49, 4).Value = '''7.56'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -1.44%  '
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(50, 2).Value = 'Algorand'
$ws.Cells.Item(50, 2).Style = "Normal"
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(50, 3).Style = "Normal"
$ws.Cells.Item(50, 4).Value = '''0.0965'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -1.02%  '
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.20%  '
$ws.Cells.Item(51, 5).Style = "Normal"
